# Bullish Engulfing - Conditions: Sharpened Condition.
# The BPCL trade (row 3) is no longer "Active"; fill in its sell-side details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E3: Active flag flips from "Yes" to "No"
$ws.Range("E3").Value = "No"

# F3: Sell Price (numeric)
$ws.Range("F3").Value = 303.6

# G3: P/L, stored as literal text (not an auto-calculated percentage)
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "+0.4%"
$ws.Range("G3").Style = "Normal"

# H3: Sell Date, stored as literal text (not an Excel date serial)
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2025-04-28"
$ws.Range("H3").Style = "Normal"

# I3: Sell Time, stored as literal text (not an Excel time serial)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "11:01:45"
$ws.Range("I3").Style = "Normal"
